$wb = $excel.ActiveWorkbook

# Sheet "2025" (rId1)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 3906.399109145206
$ws.Range("C2").Value = 48353.76274462014
$ws.Range("F2").Value = 9433.134471502228
$ws.Range("H2").Value = 2534.277928792104
$ws.Range("N2").Value = 2367.37219622158
$ws.Range("O2").Value = 1995.762462679798

# Sheet "2030" (rId2)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 6991.052031681918
$ws.Range("C2").Value = 197913.7502057619
$ws.Range("F2").Value = 16452.51445364119
$ws.Range("H2").Value = 8194.52068131253
$ws.Range("N2").Value = 7543.193583625169
$ws.Range("O2").Value = 6257.586732772244

# Sheet "2035" (rId3)
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 31236.29455387744
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("N2").Value = 12888.04225687751
$ws.Range("O2").Value = 9263.466444480218

# Sheet "2040" (rId4)
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 31236.29455387744
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("N2").Value = 14045.89200932069
$ws.Range("O2").Value = 9263.466444480218

# Sheet "2045" (rId5)
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 38906.8534480406
$ws.Range("B2").Value = 193.0947398408091
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("N2").Value = 16879.89729726143
$ws.Range("O2").Value = 10096.02314047837

# Sheet "2050" (rId6)
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 38906.8534480406
$ws.Range("B2").Value = 193.0947398408091
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("N2").Value = 16879.89729726143
$ws.Range("O2").Value = 10096.02314047837
